{"js": "// Docs - minor update\n//\n// In the \"Functional Interfaces\" section, the sentence\n//   \"...the app's deployment bar the use of some of the newer ones...\"\n// becomes\n//   \"...the app's deployment bars the use of some of the newer ones...\"\n// i.e. a single \"s\" is inserted right after the word \"bar\". The document's\n// \"_GoBack\" bookmark (which Word drops at the location of the most recent\n// edit) moves along with it, landing immediately after the newly inserted\n// \"s\" and before the following \" the use...\" text.\n\nconst body = context.document.body;\n\n// Narrow the search to the unique sentence containing the target word so we\n// don't accidentally match \"bar\" anywhere else in the document.\nconst scopedResults = body.search(\"deployment bar the use of some of the newer ones\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nscopedResults.load(\"items\");\nawait context.sync();\n\nif (scopedResults.items.length === 0) {\n  throw new Error(\"Could not locate the target sentence ('deployment bar the use of some of the newer ones').\");\n}\n\nconst sentenceRange = scopedResults.items[0];\n\n// Within that narrow range, find the exact \"bar\" token to insert after.\nconst barResults = sentenceRange.search(\"bar\", { matchCase: true });\nbarResults.load(\"items\");\nawait context.sync();\n\nif (barResults.items.length === 0) {\n  throw new Error(\"Could not locate the word 'bar' to pluralize.\");\n}\n\nconst barRange = barResults.items[0];\n\n// Insert the missing \"s\" right after \"bar\" -> \"bars\".\nbarRange.insertText(\"s\", \"After\");\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark so it sits right after the newly typed \"s\",\n// matching where Word leaves it after a real edit at that spot.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst barsResults = body.search(\"deployment bars\", { matchCase: true });\nbarsResults.load(\"items\");\nawait context.sync();\n\nif (barsResults.items.length > 0) {\n  const caret = barsResults.items[0].getRange(\"End\");\n  caret.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Docs - minor update\n#\n# In the \"Functional Interfaces\" section, the sentence\n#   \"...the app's deployment bar the use of some of the newer ones...\"\n# becomes\n#   \"...the app's deployment bars the use of some of the newer ones...\"\n# i.e. a single \"s\" is inserted right after the word \"bar\". The document's\n# \"_GoBack\" bookmark (which Word drops at the location of the most recent\n# edit) moves along with it, landing immediately after the newly inserted\n# \"s\" and before the following \" the use...\" text.\n\n$d = $word.ActiveDocument\n\n# Narrow the search to the unique sentence containing the target word so we\n# don't accidentally match \"bar\" anywhere else in the document.\n$sentence = $d.Content\n$sentence.Find.ClearFormatting()\n$sentence.Find.Text = \"deployment bar the use of some of the newer ones\"\n$sentence.Find.MatchCase = $true\n$sentence.Find.MatchWholeWord = $false\n$found = $sentence.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate the target sentence ('deployment bar the use of some of the newer ones').\"\n}\n\n# Within that sentence, find the exact \"bar\" token to insert after.\n$barRange = $sentence.Duplicate\n$barRange.Find.ClearFormatting()\n$barRange.Find.Text = \"bar\"\n$barRange.Find.MatchCase = $true\n$barRange.Find.Execute() | Out-Null\n\n# Insert the missing \"s\" right after \"bar\" -> \"bars\". InsertAfter grows the\n# calling range so $barRange.End now points right after the new \"s\".\n$barRange.InsertAfter(\"s\")\n\n# Move the \"_GoBack\" bookmark so it sits right after the newly typed \"s\",\n# matching where Word leaves it after a real edit at that spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$caret = $d.Range($barRange.End, $barRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $caret) | Out-Null\n"}
